$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill Z2:Z6 with the new "DropBag_1" value (adds a new shared string and
# references it from every row's Z column, mirroring the DropPackList header
# in Z1).
$ws.Range("Z2:Z6").Value = "DropBag_1"

# Document the intended format of the DropPackList / DropBag column with a
# cell comment on the header cell.
$commentText = "可填入英文分号间隔的奖励包(掉落包和奖励包公用)" + [char]10
$comment = $ws.Range("Z1").AddComment($commentText)

# Move the active selection the way the author left it before saving.
$ws.Range("Z6").Select() | Out-Null
